# Plantilla_Entrada_5W1H.xlsx - add new sheet "7_Categoria_Canal"
# (channel split table) at the end of the workbook and make it active.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the LAST tab ---------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "7_Categoria_Canal"

# --- Populate the data -------------------------------------------------------
$newSheet.Range("A1").Value = "Measures = Weighted R_VOL1 Vert %"
$newSheet.Range("A2").Value = "_PERIODS = MAT Jun-25\Total _PERIODS"

$newSheet.Range("A3").Value = "crch54kw - table - 05/12/2025 05:47:34 p. m."
$newSheet.Range("B3").Value = "T. Quesos Blancos+ Fundidos"
$newSheet.Range("C3").Value = "La Serenísima Clásico"
$newSheet.Range("D3").Value = "Casancrem"

$newSheet.Range("A4").Value = "Total BFPCH  "
$newSheet.Range("B4").Value = 100
$newSheet.Range("C4").Value = 100
$newSheet.Range("D4").Value = 100

$newSheet.Range("A5").Value = " Modern Trade  "
$newSheet.Range("B5").Value = 74.6
$newSheet.Range("C5").Value = 72.2
$newSheet.Range("D5").Value = 75.9

$newSheet.Range("A6").Value = " Traditional Trade  "
$newSheet.Range("B6").Value = 24.3
$newSheet.Range("C6").Value = 27.1
$newSheet.Range("D6").Value = 23.4

$newSheet.Range("A7").Value = " Other Channels  "
$newSheet.Range("B7").Value = 1.1
$newSheet.Range("C7").Value = 0.7
$newSheet.Range("D7").Value = 0.7

# --- Match the selection left behind on the new (now active) sheet ---------
$newSheet.Range("E14").Select() | Out-Null
